$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new date
$ws.Name = "PickAndPlace_PCB1_2026-01-19"

# Update the CPL coordinate values that changed in this revision
# U3 (row 2): Mid X, Ref X, Pad X
$ws.Range("D2").Value = "23.165mm"
$ws.Range("F2").Value = "23.165mm"
$ws.Range("H2").Value = "23.165mm"

# CN1 (row 3): Mid X, Ref X, Pad X
$ws.Range("D3").Value = "16.713mm"
$ws.Range("F3").Value = "16.713mm"
$ws.Range("H3").Value = "15.463mm"

# U1 (row 5): Mid X, Ref X, Pad X
$ws.Range("D5").Value = "19.863mm"
$ws.Range("F5").Value = "19.863mm"
$ws.Range("H5").Value = "27.864mm"
